$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Create Nationality", "PASSED", "chrome", "07.09.22 18:02"),
    @("Edit Nationality", "PASSED", "chrome", "07.09.22 18:29"),
    @("Delete Nationality", "PASSED", "chrome", "07.09.22 18:50"),
    @("Create Discount", "PASSED", "chrome", "07.09.22 18:39"),
    @("Edit Discount", "FAILED", "chrome", "07.09.22 18:59"),
    @("Delete Discount", "FAILED", "chrome", "07.09.22 18:52"),
    @("Create Discount", "PASSED", "chrome", "07.09.22 18:52"),
    @("Edit Discount", "PASSED", "chrome", "07.09.22 18:16"),
    @("Delete Discount", "PASSED", "chrome", "07.09.22 18:39"),
    @("Create Grade Level", "FAILED", "chrome", "07.09.22 18:55"),
    @("Edit Grade Level", "PASSED", "chrome", "07.09.22 18:20"),
    @("Delete Grade Level", "PASSED", "chrome", "07.09.22 18:43")
)

$row = 10
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row = $row + 1
}
